# Human Activity Recognition.pptx - "Linear Regression" -> "Logistic Regression"
# rename + related layout nudges on the Results slide (12) and the
# corresponding bullet-point updates on the Observations slide (13).

# PowerPoint's COM properties (Shape.Left/Top/Width/Height) are IEEE-754
# single-precision floats measured in points; converting an EMU value to
# points and back can truncate by a single EMU. EmuToPt nudges the point
# value upward (by sub-float increments) until the float32 round-trip lands
# back on the exact EMU we want, so the written XML matches byte-for-byte.
function EmuToPt {
    param([double]$emu)
    $pt = $emu / 12700.0
    for ($i = 0; $i -lt 100; $i++) {
        $f = [single]$pt
        $back = [int64]([double]$f * 12700.0)
        if ($back -eq $emu) {
            return $pt
        }
        $pt += 0.0000005
    }
    return $pt
}

# Replace a paragraph's text in one shot (selecting the paragraph's full
# character range first) so the run is rewritten as a single <a:r> with its
# original <a:rPr> preserved, instead of being split at the edit boundary.
function SetParagraphText {
    param($paragraph, [string]$newText)
    $fullRange = $paragraph.Characters(1, $paragraph.Length)
    $fullRange.Text = $newText
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 12 ("Results") - two results tables both list "Linear Regression"
# ---------------------------------------------------------------------
$s12 = $p.Slides.Item(12)

# Title placeholder gains an explicit (slightly raised) position.
$title = $s12.Shapes.Item(1)
$title.Left = EmuToPt 838200
$title.Top = EmuToPt 318472
$title.Width = EmuToPt 10515600
$title.Height = EmuToPt 1325563

# First table (shape 2): row 8 / column 1 = "Linear Regression".
$table1 = $s12.Shapes.Item(2).Table
SetParagraphText $table1.Cell(8, 1).Shape.TextFrame.TextRange.Paragraphs(1) "Logistic Regression"

# White backdrop rectangle behind the second table moves down slightly.
$rect7 = $s12.Shapes.Item(9)
$rect7.Left = EmuToPt 930987
$rect7.Top = EmuToPt 1405516
$rect7.Width = EmuToPt 9018556
$rect7.Height = EmuToPt 3233349

# Second table (shape 10) shifts up-left slightly.
$table2Shape = $s12.Shapes.Item(10)
$table2Shape.Left = EmuToPt 930987
$table2Shape.Top = EmuToPt 1388706
$table2Shape.Width = EmuToPt 8772850
$table2Shape.Height = EmuToPt 2947560

# Second table: row 8 / column 1 = "Linear Regression".
$table2 = $table2Shape.Table
SetParagraphText $table2.Cell(8, 1).Shape.TextFrame.TextRange.Paragraphs(1) "Logistic Regression"

# Four highlight rectangles around the second table nudge slightly.
$rect9 = $s12.Shapes.Item(11)
$rect9.Left = EmuToPt 2970246
$rect9.Top = EmuToPt 2081076

$rect11 = $s12.Shapes.Item(12)
$rect11.Left = EmuToPt 2970246
$rect11.Top = EmuToPt 1753901

$rect13 = $s12.Shapes.Item(13)
$rect13.Left = EmuToPt 7955384
$rect13.Top = EmuToPt 3699769

$rect14 = $s12.Shapes.Item(14)
$rect14.Left = EmuToPt 7955384
$rect14.Top = EmuToPt 2401603

# ---------------------------------------------------------------------
# Slide 13 ("Observations") - bullet text mentions of "Linear regression"
# ---------------------------------------------------------------------
$s13 = $p.Slides.Item(13)
$obsRange = $s13.Shapes.Item(2).TextFrame.TextRange

SetParagraphText $obsRange.Paragraphs(2) "Logistic regression is fast compared to others but the accuracy is relatively low."
SetParagraphText $obsRange.Paragraphs(6) "Logistic regression outperformed Ridge Regression"
